$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Insert a new column at I (pushing D.Hopkins .. Z.Ertz one column to the right)
    # so the new player "G.Dortch" sits right after "E.Benjamin".
    $ws.Columns("I:I").Insert()

    # New header cell: bold, centered, bordered - same look as the rest of row 1 (B1:S1/T1).
    $ws.Range("I1").Value = "G.Dortch"
    $ws.Range("I1").Font.Bold = $true
    $ws.Range("I1").HorizontalAlignment = -4108   # xlCenter
    $ws.Range("I1").VerticalAlignment = -4160     # xlTop
    $ws.Range("I1").Borders.LineStyle = 1         # xlContinuous

    # New data cell in row 2 matches the rest of the row ("n").
    $ws.Range("I2").Value = "n"
}
